# conditions part xpath element renamed to xpathref
#
# The "conditions" custom-xml example in the document shows:
#
#   <od:conditions xmlns:od="http://opendope.org/conditions" >
#     <od:condition id="c5">
#       <od:xpath id="x5"/>
#     </od:condition>
#     <od:condition id="c6">
#       <od:xpath id="x6"/>
#     </od:condition>
#   </od:conditions>
#
# Both <od:xpath> elements nested inside <od:condition> are renamed to
# <od:xpathref>. (A later, unrelated example re-using id="x5" - in the
# "In this example, the condition simply uses the value of ..." paragraph -
# keeps the old <od:xpath> element name, so we must not touch it.)

$d = $word.ActiveDocument

# Scope the replacement to just the <od:conditions>...</od:conditions>
# sample block so the other 16 "od:xpath" occurrences elsewhere in the
# document (and the unrelated, later re-use of id="x5") are left alone.
$scanStart = $d.Content
$null = $scanStart.Find.Execute("od:conditions xmlns:od", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$blockStart = $scanStart.Start

$scanEnd = $d.Content
$null = $scanEnd.Find.Execute("In this example", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$blockEnd = $scanEnd.Start

# First <od:xpath id="x5"/> inside <od:condition id="c5">
$work = $d.Range($blockStart, $blockEnd)
$null = $work.Find.Execute("od:xpath", $true, $false, $false, $false, $false, $true, 1, $false, "od:xpathref", 1)

# Second <od:xpath id="x6"/> inside <od:condition id="c6"> - resume searching
# right after the first replacement, still bounded by the sample block end.
$work2 = $d.Range($work.End, $blockEnd)
$null = $work2.Find.Execute("od:xpath", $true, $false, $false, $false, $false, $true, 1, $false, "od:xpathref", 1)

# --- Other textual edits in the same revision ---

# Document version bump.
$null = $d.Content.Find.Execute("document version v2", $true, $false, $false, $false, $false, $true, 1, $false, "document version v2.01", 2)

# Revised date.
$null = $d.Content.Find.Execute("10 October 2010", $true, $false, $false, $false, $false, $true, 1, $false, "14 October 2010", 2)

# Tag-length limit correction (74 -> 64 characters).
$null = $d.Content.Find.Execute("1. content in the sdt tag is minimised, which is necessary since Word restricts the tag content to 74 characters", $true, $false, $false, $false, $false, $true, 1, $false, "1. content in the sdt tag is minimised, which is necessary since Word restricts the tag content to 64 characters", 2)
